$wb = $excel.ActiveWorkbook

# Updated "想去人数" (want-to-go count) values for rows 2-16 (row 9 and 13 unchanged)
$updates = @{
    2  = 181
    3  = 224
    4  = 256
    5  = 778
    6  = 229
    7  = 5705
    8  = 20
    10 = 94
    11 = 38
    12 = 26
    14 = 177
    15 = 294
    16 = 22
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
